# Applies the "Updated cryptos list" data refresh described in the commit
# message: price (D) and 1h volume % (E) updates for most rows, plus three
# pairs of rows (19/20, 33/34, 39/40) whose coin name+link+data swapped order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.032.66'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').Value = '2.828.76'
$ws.Range('E3').Value = '  +3.61%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''354.26'
$ws.Range('E5').Value = '  +7.02%  '
$ws.Range('D6').Value = '''114.02'
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('D7').Value = '''0.550'
$ws.Range('E7').Value = '  +2.48%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '''0.604'
$ws.Range('E9').Value = '  +6.61%  '
$ws.Range('D10').Value = '''42.09'
$ws.Range('E10').Value = '  +1.95%  '
$ws.Range('D11').Value = '''0.0850'
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('D12').Value = '''20.15'
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').Value = '''0.131'
$ws.Range('E13').Value = '  +1.33%  '
$ws.Range('D14').Value = '''7.74'
$ws.Range('E14').Value = '  +2.50%  '
$ws.Range('D15').Value = '3.242.88'
$ws.Range('E15').Value = '  +2.42%  '
$ws.Range('D16').Value = '2.814.45'
$ws.Range('E16').Value = '  +1.71%  '
$ws.Range('D17').Value = '''0.899'
$ws.Range('E17').Value = '  +2.79%  '
$ws.Range('D18').Value = '51.920.17'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '''7.30'
$ws.Range('E19').Value = '  +7.37%  '
$ws.Range('B20').Value = 'ImmutableX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D20').Value = '''3.17'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('D21').Value = '''13.64'
$ws.Range('E21').Value = '  +2.52%  '
$ws.Range('D22').Value = '''0.0000100'
$ws.Range('E22').Value = '  +3.64%  '
$ws.Range('D23').Value = '''269.93'
$ws.Range('E23').Value = '  -2.98%  '
$ws.Range('D24').Value = '''69.67'
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('D25').Value = '''2.80'
$ws.Range('E25').Value = '  +6.10%  '
$ws.Range('D26').Value = '''26.75'
$ws.Range('E26').Value = '  +0.68%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('E28').Value = '  +1.82%  '
$ws.Range('D29').Value = '''2.25'
$ws.Range('E29').Value = '  +1.46%  '
$ws.Range('E30').Value = '  +0.54%  '
$ws.Range('D31').Value = '''50.79'
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('D32').Value = '''33.83'
$ws.Range('E32').Value = '  -2.63%  '
$ws.Range('B33').Value = 'VeChain'
$ws.Range('C33').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D33').Value = '''0.0452'
$ws.Range('E33').Value = '  +31.90%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '''5.87'
$ws.Range('E34').Value = '  +6.46%  '
$ws.Range('D35').Value = '''0.0833'
$ws.Range('E35').Value = '  +2.17%  '
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').Value = '''2.10'
$ws.Range('E37').Value = '  +1.24%  '
$ws.Range('D38').Value = '''18.50'
$ws.Range('E38').Value = '  -2.51%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = '''3.22'
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''4.84'
$ws.Range('E40').Value = '  -2.83%  '
$ws.Range('D41').Value = '''2.58'
$ws.Range('E41').Value = '  +7.79%  '
$ws.Range('D42').Value = '''128.25'
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('D43').Value = '''23.44'
$ws.Range('E43').Value = '  +1.74%  '
$ws.Range('E44').Value = '  +1.97%  '
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('D46').Value = '''3.36'
$ws.Range('E46').Value = '  +2.04%  '
$ws.Range('D47').Value = '2.080.44'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  +4.08%  '
$ws.Range('D49').Value = '''0.963'
$ws.Range('E49').Value = '  +11.80%  '
$ws.Range('D50').Value = '''5.70'
$ws.Range('E50').Value = '  +3.63%  '
$ws.Range('D51').Value = '''60.49'
$ws.Range('E51').Value = '  +1.66%  '
